# TP2 - PlanillaDeMetricas: fill in test-preparation timings, increment table
# (Clase Vector / Clase Matriz / Clase SEL / Paquete Test) data, and the
# execution-of-test timings. All dependent formulas (row totals, the
# "Resumen" section and the pie chart) recalc automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Preparación de la Prueba (increment 1): estimated / start / end time ---
$ws.Range("B4").Value = 0.0069444444444444441
$ws.Range("C4").Value = 0.54861111111111105
$ws.Range("D4").Value = 0.55625000000000002

# --- Preparación de la Prueba (increment 2): estimated / start / end time ---
$ws.Range("B8").Value = 0.027777777777777776
$ws.Range("C8").Value = 0.55694444444444446
$ws.Range("D8").Value = 0.58333333333333337

# --- Desarrollo y correctivos: increment rows 13-16 ---
$ws.Range("C13").Value = "Clase Vector"
$ws.Range("F13").Value = 150
$ws.Range("G13").Value = 0.027777777777777776
$ws.Range("H13").Value = 0.79166666666666663
$ws.Range("I13").Value = 0.81597222222222221
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 168

$ws.Range("C14").Value = "Clase Matriz"
$ws.Range("F14").Value = 200
$ws.Range("G14").Value = 0.027777777777777776
$ws.Range("H14").Value = 0.82291666666666663
$ws.Range("I14").Value = 0.86805555555555547
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.0069444444444444441
$ws.Range("M14").Value = 265

$ws.Range("C15").Value = "Clase SEL"
$ws.Range("F15").Value = 100
$ws.Range("G15").Value = 0.017361111111111112
$ws.Range("H15").Value = 0.59027777777777779
$ws.Range("I15").Value = 0.60416666666666663
$ws.Range("K15").Value = 1
$ws.Range("L15").Value = 0.003472222222222222
$ws.Range("M15").Value = 134

$ws.Range("C16").Value = "Paquete Test"
$ws.Range("F16").Value = 40
$ws.Range("G16").Value = 0.0069444444444444441
$ws.Range("H16").Value = 0.60416666666666663
$ws.Range("I16").Value = 0.60972222222222217
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 48

# Row 15 no longer needs the extra-tall custom row height now that it holds
# a short class name instead of wrapped placeholder text.
$ws.Rows.Item(15).AutoFit()

# --- Ejecución de la Prueba: estimated / start / end time ---
$ws.Range("B21").Value = 0.017361111111111112
$ws.Range("C21").Value = 0.61111111111111105
$ws.Range("D21").Value = 0.625

# Leave the cursor where the author left it after finishing data entry.
$ws.Range("C1:N1").Select()

$wb.Save()
